# FIX #12091 TIME 0:25 update template with [attachmentRecipient. intead of [recipient.
#
# Applies the changes described by the diff against
# modules/templates/templates/styles/standard_nosign.docx :
#   1. Bump the left cell margin of both tables from 128 -> 133 dxa
#      (6.4pt -> 6.65pt).
#   2. Rename the merge field placeholder "[recipient." to
#      "[attachmentRecipient." everywhere it appears (postal address
#      block, salutation line and closing line).
#   3. Refresh the cached value of the "Maarch-les-Bains, le <TIME>"
#      date field from 09/12/2019 to 02/01/2020.

$d = $word.ActiveDocument

# --- 1. Table cell left padding: 128 dxa (6.4pt) -> 133 dxa (6.65pt) ---
foreach ($t in $d.Tables) {
    $t.LeftPadding = 6.65
}

# --- 2. [recipient. -> [attachmentRecipient. ------------------------
# A single Find/Replace over the whole story also normalises the
# surrounding runs (the source document had the merge field split
# across several runs - e.g. "recipient" / ".c" / "ivility" - which
# collapse into one run per replaced field, matching the target XML).
$d.Content.Find.Execute("recipient", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "attachmentRecipient", 2) | Out-Null

# --- 3. Refresh the cached TIME field result ------------------------
$d.Content.Find.Execute("09/12/2019", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "02/01/2020", 2) | Out-Null

# --- 4. Nudge the thin horizontal-line autoshape in the header ------
# (its cached extent is recomputed by Word whenever the header gets
# touched; reassert it explicitly so the anchor extent matches).
foreach ($sec in $d.Sections) {
    foreach ($hf in $sec.Headers) {
        if ($hf.Exists) {
            foreach ($s in $hf.Shapes) {
                if ($s.Type -eq 1) {
                    $s.Width = 543.89999
                    $s.Height = 0.35
                }
            }
        }
    }
}
